$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "A"="ECs"; "B"="Fn1"; "C"="Itgb8"; "D"="ECs"; "E"=3; "F"=1; "G"=21.84976866666667; "H"=65.549306; "I"=0.05020018890879543; "J"=0.05020018890879543; "K"=1; "L"=0.3333333333333333; "M"=0.1126243333333333; "N"=0.337873; "O"=0.01082936903163217; "P"=0.01082936903163217; "Q"=2.460815629570889; "R"=22.147340666138; "S"=0.0005436363711509938; "T"=0.0005436363711509938 }
  3 = @{ "A"="ECs"; "B"="Fn1"; "C"="Itgb8"; "D"="FAPs"; "E"=3; "F"=1; "G"=21.84976866666667; "H"=65.549306; "I"=0.05020018890879543; "J"=0.05020018890879543; "K"=3; "L"=1; "M"=3.323421; "N"=9.970263; "O"=0.3195628457125252; "P"=0.3195628457125252; "Q"=72.61598003194199; "R"=653.543820287478; "S"=0.01604211522300101; "T"=0.01604211522300101 }
  4 = @{ "A"="ECs"; "B"="Fn1"; "C"="Itgb8"; "D"="sCs"; "E"=3; "F"=1; "G"=21.84976866666667; "H"=65.549306; "I"=0.05020018890879543; "J"=0.05020018890879543; "K"=3; "L"=1; "M"=6.963852666666667; "N"=20.891558; "O"=0.6696077852558425; "P"=0.6696077852558425; "Q"=152.1585697954164; "R"=1369.427128158748; "S"=0.03361443731464342; "T"=0.03361443731464342 }
  5 = @{ "A"="FAPs"; "B"="Fn1"; "C"="Itgb8"; "D"="ECs"; "E"=3; "F"=1; "G"=385.0524703333334; "H"=1155.157411; "I"=0.8846641374295412; "J"=0.8846641374295412; "K"=1; "L"=0.3333333333333333; "M"=0.1126243333333333; "N"=0.337873; "O"=0.01082936903163217; "P"=0.01082936903163217; "Q"=43.36627776964477; "R"=390.296499926803; "S"=0.009580354413275058; "T"=0.009580354413275058 }
  6 = @{ "A"="FAPs"; "B"="Fn1"; "C"="Itgb8"; "D"="FAPs"; "E"=3; "F"=1; "G"=385.0524703333334; "H"=1155.157411; "I"=0.8846641374295412; "J"=0.8846641374295412; "K"=3; "L"=1; "M"=3.323421; "N"=9.970263; "O"=0.3195628457125252; "P"=0.3195628457125252; "Q"=1279.691466007677; "R"=11517.22319406909; "S"=0.2827057892568007; "T"=0.2827057892568007 }
  7 = @{ "A"="FAPs"; "B"="Fn1"; "C"="Itgb8"; "D"="sCs"; "E"=3; "F"=1; "G"=385.0524703333334; "H"=1155.157411; "I"=0.8846641374295412; "J"=0.8846641374295412; "K"=3; "L"=1; "M"=6.963852666666667; "N"=20.891558; "O"=0.6696077852558425; "P"=0.6696077852558425; "Q"=2681.448672337371; "R"=24133.03805103634; "S"=0.5923779937594654; "T"=0.5923779937594654 }
  8 = @{ "A"="sCs"; "B"="Fn1"; "C"="Itgb8"; "D"="ECs"; "E"=3; "F"=1; "G"=28.350479; "H"=85.05143699999999; "I"=0.06513567366166337; "J"=0.06513567366166337; "K"=1; "L"=0.3333333333333333; "M"=0.1126243333333333; "N"=0.337873; "O"=0.01082936903163217; "P"=0.01082936903163217; "Q"=3.192953797055666; "R"=28.736584173501; "S"=0.0007053782472061164; "T"=0.0007053782472061164 }
  9 = @{ "A"="sCs"; "B"="Fn1"; "C"="Itgb8"; "D"="FAPs"; "E"=3; "F"=1; "G"=28.350479; "H"=85.05143699999999; "I"=0.06513567366166337; "J"=0.06513567366166337; "K"=3; "L"=1; "M"=3.323421; "N"=9.970263; "O"=0.3195628457125252; "P"=0.3195628457125252; "Q"=94.22057726865899; "R"=847.9851954179309; "S"=0.02081494123272352; "T"=0.02081494123272352 }
  10 = @{ "A"="sCs"; "B"="Fn1"; "C"="Itgb8"; "D"="sCs"; "E"=3; "F"=1; "G"=28.350479; "H"=85.05143699999999; "I"=0.06513567366166337; "J"=0.06513567366166337; "K"=3; "L"=1; "M"=6.963852666666667; "N"=20.891558; "O"=0.6696077852558425; "P"=0.6696077852558425; "Q"=197.4285587854273; "R"=1776.857029068846; "S"=0.04361535418173372; "T"=0.04361535418173372 }
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in ($data.Keys | Sort-Object)) {
    $rowVals = $data[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $rowVals[$c]
    }
}
